# Nomina quincenal - interface improvements:
#  - Insert a new column "Empleado por contrato" right after "Empleado Fijo"
#  - Insert 4 new discount columns (Seguro Social, Seguro Educativo, ISL,
#    Total Descuentos) right before "Numero de Cuenta"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert new column G ("Empleado por contrato") ---
# Column G currently holds "Salario Base"; shift it (and everything after)
# one column to the right, then fill the new G column.
$ws.Columns.Item(7).Insert()

$ws.Range("G1").Value = "Empleado por contrato"

$lastRow = 7
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = "No"
}

# --- Step 2: insert 4 new columns for tax/insurance deductions ---
# After the first insertion, "Numero de Cuenta" now lives in column R.
# Insert four columns right before it (R:U) to make room for:
#   Seguro Social (3%), Seguro Educativo (5%), ISL, Total Descuentos
$ws.Range("R1:U1").EntireColumn.Insert()

$ws.Range("R1").Value = "Seguro Social (3%)"
$ws.Range("S1").Value = "Seguro Educativo (5%)"
$ws.Range("T1").Value = "ISL"
$ws.Range("U1").Value = "Total Descuentos"

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 18).Value = 0  # R - Seguro Social (3%)
    $ws.Cells.Item($r, 19).Value = 0  # S - Seguro Educativo (5%)
    $ws.Cells.Item($r, 20).Value = 0  # T - ISL
    $ws.Cells.Item($r, 21).Value = 0  # U - Total Descuentos
}

Write-Output "Edit complete"
